$wb = $excel.ActiveWorkbook

# --- Sheet "All Orders": mark order row 10 as CANCELLED with a cancel reason ---
$wsOrders = $wb.Worksheets.Item("All Orders")
$wsOrders.Range("H10").Value = "CANCELLED"
$wsOrders.Range("M10").Value = "test order"

# --- Sheet "Daily Summary": update 2026-01-13 totals to reflect the cancellation ---
$wsSummary = $wb.Worksheets.Item("Daily Summary")
$wsSummary.Range("D4").Value = 7
$wsSummary.Range("E4").Value = 110
$wsSummary.Range("G4").Value = 110
